$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "92.362.82"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "3.102.32"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").Value = "'240.68"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").Value = "'615.25"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("E7").Value = "  -6.08%  "
$ws.Range("D8").Value = "'0.394"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.35%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "3.102.34"
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("E11").Value = "  -4.69%  "
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("B14").Value = "WrappedBTC"
$ws.Range("C14").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D14").Value = "92.114.68"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'34.38"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.29%  "
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "3.679.87"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").Value = "3.093.49"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("E19").Value = "  -2.70%  "
$ws.Range("E20").Value = "  -2.66%  "
$ws.Range("D21").Value = "'5.82"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("D22").Value = "'9.38"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("D23").Value = "'447.84"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("E24").Value = "  -5.49%  "
$ws.Range("E25").Value = "  -1.97%  "
$ws.Range("D26").Value = "'87.22"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.24%  "
$ws.Range("D27").Value = "'11.76"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("D28").Value = "3.264.35"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").ClearFormats()
$ws.Range("D30").Value = "'0.138"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +7.88%  "
$ws.Range("E31").Value = "  -3.56%  "
$ws.Range("D32").Value = "'0.167"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.03%  "
$ws.Range("D33").Value = "'9.23"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("D34").Value = "'0.997"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.28%  "
$ws.Range("D35").Value = "'8.03"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("E36").Value = "  -5.75%  "
$ws.Range("D37").Value = "'4.23"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.44%  "
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("D39").Value = "'1.93"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("E40").Value = "  -0.85%  "
$ws.Range("D41").Value = "'480.06"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.45%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "'0.434"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'3.45"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.41%  "
$ws.Range("D44").Value = "'22.88"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.73%  "
$ws.Range("D46").Value = "'160.15"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.71%  "
$ws.Range("E47").Value = "  -2.76%  "
$ws.Range("D48").Value = "'0.693"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.85%  "
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("D50").Value = "'0.0332"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.94%  "
$ws.Range("D51").Value = "'44.13"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.57%  "
